$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3299.6
$ws.Cells.Item(64, 9).Value = 2999
$ws.Cells.Item(64, 11).Value = 2999
$ws.Cells.Item(64, 13).Value = -2751

$ws.Cells.Item(67, 8).Value = 3299.6
$ws.Cells.Item(67, 9).Value = 2999
$ws.Cells.Item(67, 11).Value = 2999
$ws.Cells.Item(67, 13).Value = -2141

$ws.Cells.Item(86, 8).Value = 2172574.5
$ws.Cells.Item(86, 9).Value = 12006.889
$ws.Cells.Item(86, 10).Value = 6061596
$ws.Cells.Item(86, 11).Value = 12006.889
$ws.Cells.Item(86, 12).Value = 6061596
$ws.Cells.Item(86, 13).Value = -10883.889
$ws.Cells.Item(86, 14).Value = -6063842

$ws.Cells.Item(89, 8).Value = 2172574.5
$ws.Cells.Item(89, 9).Value = 12006.889
$ws.Cells.Item(89, 10).Value = 6061596
$ws.Cells.Item(89, 11).Value = 60034.44499999999
$ws.Cells.Item(89, 12).Value = 30307980
$ws.Cells.Item(89, 13).Value = -54418.44499999999
$ws.Cells.Item(89, 14).Value = -30319212

$ws.Cells.Item(100, 8).Value = 12789.789
$ws.Cells.Item(100, 9).Value = 26488.75
$ws.Cells.Item(100, 10).Value = 2826.9092
$ws.Cells.Item(100, 11).Value = 26488.75
$ws.Cells.Item(100, 12).Value = 2826.9092
$ws.Cells.Item(100, 13).Value = -25947.75
$ws.Cells.Item(100, 14).Value = -3908.9092

$ws.Cells.Item(137, 8).Value = 1470.8182
$ws.Cells.Item(137, 9).Value = 1185.1428
$ws.Cells.Item(137, 10).Value = 3070.6
$ws.Cells.Item(137, 11).Value = 3555.4284
$ws.Cells.Item(137, 12).Value = 9211.799999999999
$ws.Cells.Item(137, 13).Value = -1005.4284
$ws.Cells.Item(137, 14).Value = -14311.8

$ws.Cells.Item(138, 8).Value = 2591.2122
$ws.Cells.Item(138, 9).Value = 1048.6923
$ws.Cells.Item(138, 10).Value = 3593.85
$ws.Cells.Item(138, 11).Value = 3146.0769
$ws.Cells.Item(138, 12).Value = 10781.55
$ws.Cells.Item(138, 13).Value = 1993.9231
$ws.Cells.Item(138, 14).Value = -21061.55

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 13003461
$ws.Cells.Item(32, 9).Value = 15163129
$ws.Cells.Item(32, 10).Value = 45454.363
$ws.Cells.Item(32, 11).Value = 15163129
$ws.Cells.Item(32, 12).Value = 45454.363
$ws.Cells.Item(32, 13).Value = -15162842
$ws.Cells.Item(32, 14).Value = -46028.363

$ws.Cells.Item(45, 8).Value = 25808.44
$ws.Cells.Item(45, 9).Value = 29912.4
$ws.Cells.Item(45, 10).Value = 1868.6666
$ws.Cells.Item(45, 11).Value = 29912.4
$ws.Cells.Item(45, 12).Value = 1868.6666
$ws.Cells.Item(45, 13).Value = -29535.4
$ws.Cells.Item(45, 14).Value = -2622.6666

$ws.Cells.Item(61, 8).Value = 1930.1
$ws.Cells.Item(61, 9).Value = 1962.5927
$ws.Cells.Item(61, 10).Value = 1862.6154
$ws.Cells.Item(61, 11).Value = 1962.5927
$ws.Cells.Item(61, 12).Value = 1862.6154
$ws.Cells.Item(61, 13).Value = -1750.5927
$ws.Cells.Item(61, 14).Value = -2286.6154

$ws.Cells.Item(88, 8).Value = 2562.75
$ws.Cells.Item(88, 9).Value = 2549.875
$ws.Cells.Item(88, 10).Value = 2571.3333
$ws.Cells.Item(88, 11).Value = 2549.875
$ws.Cells.Item(88, 12).Value = 2571.3333
$ws.Cells.Item(88, 13).Value = -2143.875
$ws.Cells.Item(88, 14).Value = -3383.3333

$ws.Cells.Item(91, 8).Value = 2562.75
$ws.Cells.Item(91, 9).Value = 2549.875
$ws.Cells.Item(91, 10).Value = 2571.3333
$ws.Cells.Item(91, 11).Value = 2549.875
$ws.Cells.Item(91, 12).Value = 2571.3333
$ws.Cells.Item(91, 13).Value = -1145.875
$ws.Cells.Item(91, 14).Value = -5379.3333

$ws.Cells.Item(98, 8).Value = 21000
$ws.Cells.Item(98, 10).Value = 21000
$ws.Cells.Item(98, 12).Value = 21000
$ws.Cells.Item(98, 14).Value = -26990

$ws.Cells.Item(122, 8).Value = 5995.077
$ws.Cells.Item(122, 9).Value = 7523.6
$ws.Cells.Item(122, 11).Value = 22570.8
$ws.Cells.Item(122, 13).Value = -20120.8

$ws.Cells.Item(136, 8).Value = 1930.1
$ws.Cells.Item(136, 9).Value = 1962.5927
$ws.Cells.Item(136, 10).Value = 1862.6154
$ws.Cells.Item(136, 11).Value = 5887.7781
$ws.Cells.Item(136, 12).Value = 5587.8462
$ws.Cells.Item(136, 13).Value = -3337.7781
$ws.Cells.Item(136, 14).Value = -10687.8462

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1573.8235
$ws.Cells.Item(99, 9).Value = 1181.2285
$ws.Cells.Item(99, 10).Value = 2432.625
$ws.Cells.Item(99, 11).Value = 1181.2285
$ws.Cells.Item(99, 12).Value = 2432.625
$ws.Cells.Item(99, 13).Value = 316.7715000000001
$ws.Cells.Item(99, 14).Value = -5428.625

$ws.Cells.Item(107, 8).Value = 83334080
$ws.Cells.Item(107, 9).Value = 100000600
$ws.Cells.Item(107, 11).Value = 100000600
$ws.Cells.Item(107, 13).Value = -99998680

$ws.Cells.Item(108, 8).Value = 33907
$ws.Cells.Item(108, 10).Value = 33907
$ws.Cells.Item(108, 12).Value = 33907
$ws.Cells.Item(108, 14).Value = -41587

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1715.2727
$ws.Cells.Item(31, 9).Value = 992.8182
$ws.Cells.Item(31, 10).Value = 2293.2363
$ws.Cells.Item(31, 11).Value = 992.8182
$ws.Cells.Item(31, 12).Value = 2293.2363
$ws.Cells.Item(31, 13).Value = -697.8182
$ws.Cells.Item(31, 14).Value = -2883.2363

$ws.Cells.Item(34, 8).Value = 1715.2727
$ws.Cells.Item(34, 9).Value = 992.8182
$ws.Cells.Item(34, 10).Value = 2293.2363
$ws.Cells.Item(34, 11).Value = 992.8182
$ws.Cells.Item(34, 12).Value = 2293.2363
$ws.Cells.Item(34, 13).Value = -790.8182
$ws.Cells.Item(34, 14).Value = -2697.2363

$ws.Cells.Item(99, 8).Value = 76933200
$ws.Cells.Item(99, 9).Value = 100011860
$ws.Cells.Item(99, 11).Value = 100011860
$ws.Cells.Item(99, 13).Value = -100010362

$ws.Cells.Item(126, 8).Value = 76933200
$ws.Cells.Item(126, 9).Value = 100011860
$ws.Cells.Item(126, 11).Value = 300035580
$ws.Cells.Item(126, 13).Value = -300033110

$ws.Cells.Item(134, 8).Value = 1669.3778
$ws.Cells.Item(134, 9).Value = 1665.4117
$ws.Cells.Item(134, 10).Value = 1681.6364
$ws.Cells.Item(134, 11).Value = 4996.2351
$ws.Cells.Item(134, 12).Value = 5044.9092
$ws.Cells.Item(134, 13).Value = -2461.2351
$ws.Cells.Item(134, 14).Value = -10114.9092

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 43483496
$ws.Cells.Item(107, 9).Value = 235
$ws.Cells.Item(107, 10).Value = 52637864
$ws.Cells.Item(107, 11).Value = 705
$ws.Cells.Item(107, 12).Value = 157913592
$ws.Cells.Item(107, 13).Value = 1215
$ws.Cells.Item(107, 14).Value = -157917432

$ws.Cells.Item(131, 8).Value = 754.4400000000001
$ws.Cells.Item(131, 9).Value = 480
$ws.Cells.Item(131, 10).Value = 784.93335
$ws.Cells.Item(131, 11).Value = 1440
$ws.Cells.Item(131, 12).Value = 2354.80005
$ws.Cells.Item(131, 13).Value = 3600
$ws.Cells.Item(131, 14).Value = -12434.80005

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(64, 8).Value = 14298.2
$ws.Cells.Item(64, 9).Value = 8246
$ws.Cells.Item(64, 10).Value = 18333
$ws.Cells.Item(64, 11).Value = 8246
$ws.Cells.Item(64, 12).Value = 18333
$ws.Cells.Item(64, 13).Value = -7998
$ws.Cells.Item(64, 14).Value = -18829

$ws.Cells.Item(67, 8).Value = 14298.2
$ws.Cells.Item(67, 9).Value = 8246
$ws.Cells.Item(67, 10).Value = 18333
$ws.Cells.Item(67, 11).Value = 8246
$ws.Cells.Item(67, 12).Value = 18333
$ws.Cells.Item(67, 13).Value = -7388
$ws.Cells.Item(67, 14).Value = -20049

$ws.Cells.Item(102, 8).Value = 1507
$ws.Cells.Item(102, 9).Value = 1500
$ws.Cells.Item(102, 11).Value = 1500
$ws.Cells.Item(102, 13).Value = 122

$ws.Cells.Item(122, 8).Value = 15631691
$ws.Cells.Item(122, 9).Value = 23817328
$ws.Cells.Item(122, 10).Value = 4564.273
$ws.Cells.Item(122, 11).Value = 71451984
$ws.Cells.Item(122, 12).Value = 13692.819
$ws.Cells.Item(122, 13).Value = -71449534
$ws.Cells.Item(122, 14).Value = -18592.819

$ws.Cells.Item(126, 8).Value = 4112.7334
$ws.Cells.Item(126, 9).Value = 3171.4285
$ws.Cells.Item(126, 10).Value = 4936.375
$ws.Cells.Item(126, 11).Value = 9514.2855
$ws.Cells.Item(126, 12).Value = 14809.125
$ws.Cells.Item(126, 13).Value = -7044.2855
$ws.Cells.Item(126, 14).Value = -19749.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1993.0741
$ws.Cells.Item(7, 9).Value = 1743.238
$ws.Cells.Item(7, 10).Value = 2867.5
$ws.Cells.Item(7, 11).Value = 1743.238
$ws.Cells.Item(7, 12).Value = 2867.5
$ws.Cells.Item(7, 13).Value = -1631.238
$ws.Cells.Item(7, 14).Value = -3091.5

$ws.Cells.Item(22, 8).Value = 2534165.5
$ws.Cells.Item(22, 9).Value = 4219942.5
$ws.Cells.Item(22, 11).Value = 4219942.5
$ws.Cells.Item(22, 13).Value = -4219647.5

$ws.Cells.Item(27, 8).Value = 2534165.5
$ws.Cells.Item(27, 9).Value = 4219942.5
$ws.Cells.Item(27, 11).Value = 4219942.5
$ws.Cells.Item(27, 13).Value = -4219835.5

$ws.Cells.Item(122, 8).Value = 16232.389
$ws.Cells.Item(122, 9).Value = 16599
$ws.Cells.Item(122, 10).Value = 10000
$ws.Cells.Item(122, 11).Value = 49797
$ws.Cells.Item(122, 12).Value = 30000
$ws.Cells.Item(122, 13).Value = -47347
$ws.Cells.Item(122, 14).Value = -34900

$ws.Cells.Item(126, 8).Value = 1993.0741
$ws.Cells.Item(126, 9).Value = 1743.238
$ws.Cells.Item(126, 10).Value = 2867.5
$ws.Cells.Item(126, 11).Value = 5229.714
$ws.Cells.Item(126, 12).Value = 8602.5
$ws.Cells.Item(126, 13).Value = -2759.714
$ws.Cells.Item(126, 14).Value = -13542.5
